$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header section
$ws.Range("G4").Value = 44084.78196113411
$ws.Range("G7").Value = "Los marineros"
$ws.Range("F10").Value = "test"

# Row 18 - piqlConnect (only piqlFilm)
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 1500
$ws.Range("H18").Value = 1500

# Row 19 - Digital (GB)
$ws.Range("F19").Value = 450
$ws.Range("G19").Value = 15
$ws.Range("H19").Value = 6750

# Row 20 - Visual (pages)
# E20 is stored as text ("2") even though the cell's number format is
# numeric; force text entry, then restore the original number format so
# the cell style index is unaffected.
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2"
$ws.Range("E20").NumberFormat = "0.00"
$ws.Range("F20").Value = 780000
$ws.Range("G20").Value = 0.016
$ws.Range("H20").Value = 12480

# Row 21 - Online Storage (GB) - piqlConnect, cleared
$ws.Range("F21").ClearContents()
$ws.Range("G21").ClearContents()
$ws.Range("H21").ClearContents()

# Row 22 - Online Storage (GB) payment info, cleared
$ws.Range("E22").ClearContents()
$ws.Range("F22").ClearContents()
$ws.Range("G22").ClearContents()
$ws.Range("H22").ClearContents()

# Row 29 - piqlReader, cleared
$ws.Range("F29").ClearContents()
$ws.Range("G29").ClearContents()
$ws.Range("H29").ClearContents()

# Row 30 - Installation and training, cleared
$ws.Range("F30").ClearContents()
$ws.Range("G30").ClearContents()
$ws.Range("H30").ClearContents()

# Row 31 - Service agreement, cleared
$ws.Range("E31").ClearContents()
$ws.Range("F31").ClearContents()
$ws.Range("G31").ClearContents()
$ws.Range("H31").ClearContents()

# Row 32 - Shipment cost
$ws.Range("E32").Value = 10
$ws.Range("H32").Value = 300

# Row 33 - TOTAL
$ws.Range("H33").Value = 13980

# Row 34 - Total to pay from the second term
$ws.Range("H34").Value = 0
